$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp label in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 18:22"

# --- Swap country labels where new data overtook the previous rank ---
# Egipto overtakes Sudafrica
$ws.Range("A52").Value = "Egipto"
$ws.Range("A53").Value = "Sudafrica"

# Somalia overtakes Georgia
$ws.Range("A108").Value = "Somalia"
$ws.Range("A109").Value = "Georgia"

# --- Refreshed numeric data (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1013557
$ws.Range("C4").Value = 3201
$ws.Range("D4").Value = 139481
$ws.Range("E4").Value = 817027

# Italia (row 6)
$ws.Range("B6").Value = 201505
$ws.Range("C6").Value = 2091
$ws.Range("D6").Value = 68941
$ws.Range("E6").Value = 105205
$ws.Range("F6").Value = 1863
$ws.Range("G6").Value = 382
$ws.Range("H6").Value = 27359

# Turquia (row 10)
$ws.Range("B10").Value = 114653
$ws.Range("C10").Value = 2392
$ws.Range("E10").Value = 77870
$ws.Range("F10").Value = 1621
$ws.Range("G10").Value = 92
$ws.Range("H10").Value = 2992

# Rumania (row 36)
$ws.Range("E36").Value = 7549
$ws.Range("G36").Value = 22
$ws.Range("H36").Value = 663

# Chequia (row 45)
$ws.Range("B45").Value = 7486
$ws.Range("C45").Value = 41
$ws.Range("D45").Value = 2942
$ws.Range("E45").Value = 4319
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 225

# Republica Dominicana (row 48)
$ws.Range("B48").Value = 6416
$ws.Range("C48").Value = 123
$ws.Range("D48").Value = 1165
$ws.Range("E48").Value = 4965
$ws.Range("G48").Value = 4
$ws.Range("H48").Value = 286

# Egipto, now row 52
$ws.Range("B52").Value = 5042
$ws.Range("C52").Value = 260
$ws.Range("D52").Value = 1304
$ws.Range("E52").Value = 3379
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 22
$ws.Range("H52").Value = 359

# Sudafrica, now row 53
$ws.Range("B53").Value = 4793
$ws.Range("D53").Value = 1473
$ws.Range("E53").Value = 3230
$ws.Range("F53").Value = 36
$ws.Range("H53").Value = 90

# Somalia, now row 108
$ws.Range("B108").Value = 528
$ws.Range("C108").Value = 48
$ws.Range("D108").Value = 19
$ws.Range("E108").Value = 481
$ws.Range("F108").Value = 2
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = 28

# Georgia, now row 109
$ws.Range("B109").Value = 511
$ws.Range("C109").Value = 14
$ws.Range("D109").Value = 156
$ws.Range("E109").Value = 349
$ws.Range("F109").Value = 6
$ws.Range("H109").Value = 6
